$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# "Step 1: cyclic_check()" bullet list - simplify the instructions.
#
# Before (one bullet, with a manual line break in the middle):
#   "Run DFS and BFS at v = 0, then print results saved in the graph
#    structure.<break>Don't use print_DFS() and print_BFS(). "
#
# After (split into two separate bullets):
#   "Run DFS and BFS at v = 0."
#   "Print results saved in the graph structure. Use print_DFS() and
#    print_BFS(). "
#
# Replacing the line break with a paragraph mark ("^p") turns the
# second half of the old bullet into its own list paragraph (it
# inherits the same list/paragraph formatting automatically), and we
# reword the lead-in text of that new paragraph at the same time.
# ---------------------------------------------------------------------
$oldText = ", then print results saved in the graph structure." + [char]11 + "Don't use "
$newText = "." + "^p" + "Print results saved in the graph structure. Use "
$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null

# ---------------------------------------------------------------------
# The hidden "_GoBack" bookmark (Word's "last edit position" marker)
# moves from the end of the document section (right after "...graphx
# .exe provided.") to sit between "an" and "d" of "print_DFS() and
# print_BFS()" - i.e. right where the author's cursor ended up after
# the edit above. Bookmarks.Add re-uses the existing "_GoBack" name,
# which relocates it (Word keeps bookmark names unique), so the old
# occurrence disappears on its own.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("print_DFS() and print_BFS()") | Out-Null
$splitAt = $rng.Start + "print_DFS() an".Length
$d.Bookmarks.Add("_GoBack", $d.Range($splitAt, $splitAt)) | Out-Null
